# The reaction-sensitivity workflow was re-run with an updated cutoff
# window: for both the NBR and BAR sheets, the first 4 data rows
# (Cutoff = 0..3) are dropped and the remaining rows shift up, keeping
# their original "Cutoff_index"/"Reaction_number" pairing but being
# renumbered sequentially in column A (0..14).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Capture the original B (Cutoff) / C (Reaction_number) values for the
    # 19 data rows (rows 2..20) before any mutation.
    $origB = @{}
    $origC = @{}
    for ($r = 2; $r -le 20; $r++) {
        $origB[$r] = $ws.Cells.Item($r, 2).Value2
        $origC[$r] = $ws.Cells.Item($r, 3).Value2
    }

    # Drop the last 4 data rows (rows 17..20) -- this shrinks the used
    # range from A1:C20 down to A1:C16, matching the new dimension.
    $ws.Range("A17:C20").EntireRow.Delete()

    # The surviving rows (new rows 2..16) take on the B/C values that used
    # to live 4 rows further down (old rows 6..20), while column A keeps a
    # fresh sequential index starting at 0.
    for ($r = 2; $r -le 16; $r++) {
        $srcRow = $r + 4
        $ws.Cells.Item($r, 2).Value = $origB[$srcRow]
        $ws.Cells.Item($r, 3).Value = $origC[$srcRow]
    }
}
